# Insert a new data row before the current row 168 in the Berenjena (Vega
# Monumental Concepción) price list, shifting all subsequent rows down by
# one. The inserted row carries a new weekly price observation; all rows
# that used to be 168..192 become 169..193 (their contents are unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 168:192 down to 169:193, opening up a blank row 168.
$ws.Rows.Item(168).Insert()

# Populate the newly inserted row 168 with the new observation.
$ws.Cells.Item(168, 1).Value  = 11
$ws.Cells.Item(168, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(168, 3).Value  = "Bíobío"
$ws.Cells.Item(168, 4).Value  = 45209
$ws.Cells.Item(168, 5).Value  = 8
$ws.Cells.Item(168, 6).Value  = 100112001
$ws.Cells.Item(168, 7).Value  = "Berenjena"
$ws.Cells.Item(168, 8).Value  = "Sin especificar"
$ws.Cells.Item(168, 9).Value  = "Primera"
$ws.Cells.Item(168, 10).Value = 80
$ws.Cells.Item(168, 11).Value = 9000
$ws.Cells.Item(168, 12).Value = 9000
$ws.Cells.Item(168, 13).Value = 9000
$ws.Cells.Item(168, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(168, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(168, 16).Value = 180
$ws.Cells.Item(168, 17).Value = 50
$ws.Cells.Item(168, 18).Value = "Hortaliza"
